$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update values for B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 are deleted (cleared), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -5.5736298643209068
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -5.8848721239177095

# Row 3: update values for B3:E3
$ws.Range("B3").Value = -6.240402481620599
$ws.Range("C3").Value = -1.8689392852975644
$ws.Range("D3").Value = -8.9813355715450243
$ws.Range("E3").Value = 9.1157060016048774

# Update selection to match new range
$ws.Range("B1:E3").Select()
